# Atualização de bases das ligas, do dia: 28-05-2024 às 07:50
#
# This edit corrects mismatched rows in the "Colombia Primera A" sheet: the
# match statistics (columns B, C and E through AD) had been written against
# the wrong fixture id/date rows. We restore the correct pairing by
# re-distributing the row contents among the affected rows, while leaving
# column A (the sequential row id) and column D (the match date, which is
# identical across every row in each affected group) untouched.
#
# Row 12 <- old Row 13, Row 13 <- old Row 12
# Row 425 <- old Row 426, Row 426 <- old Row 425
# Row 427 <- old Row 431, Row 428 <- old Row 427, Row 429 <- old Row 428,
# Row 430 <- old Row 429, Row 431 <- old Row 430

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each group lists the rows, in order, that form a closed cycle of moves:
# contents of group[i] move into group[i-1], and group[0] receives the
# contents that were in the last element of the group (wrap-around).
$groups = @(
    ,@(13, 12)
    ,@(426, 425)
    ,@(428, 427, 431, 430, 429)
)

foreach ($group in $groups) {
    $count = $group.Length

    # Snapshot the current B:AD values for every row in this group first,
    # since the rows depend cyclically on each other.
    $snapshots = @{}
    foreach ($row in $group) {
        $snapshots[$row] = $ws.Range("B$row`:AD$row").Value()
    }

    # destRow (group[i]) receives the snapshot captured from srcRow
    # (group[(i+1) mod count]).
    for ($i = 0; $i -lt $count; $i++) {
        $destRow = $group[$i]
        $srcRow = $group[($i + 1) % $count]
        $ws.Range("B$destRow`:AD$destRow").Value = $snapshots[$srcRow]
    }
}
